$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptos list (prices + 1h volume deltas), matching the
# "Updated cryptos list ... with GitHub Actions" commit. A few rows also
# swapped which coin occupies which rank (34/35, 38/39, 45/47).
#
# Price cells that are plain decimal numbers (e.g. "186.44") get forced
# back to text via NumberFormat "@" so Excel doesn't silently convert them
# to numeric cells; Style is then reset to "Normal" so no stray number
# format sticks on the cell. Prices already containing two dots
# (e.g. "67.414.97") can never parse as a number, so they're set directly.
$ws.Range("D2").Value = '67.414.97'
$ws.Range("E2").Value = '  -1.12%  '
$ws.Range("D3").Value = '3.312.89'
$ws.Range("E3").Value = '  +1.16%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '186.44'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.72%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '577.52'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.83%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("E8").Value = '  -0.31%  '
$ws.Range("E9").Value = '  -0.73%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.65'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.82%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.410'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.07%  '
$ws.Range("D12").Value = '3.892.25'
$ws.Range("E12").Value = '  +1.26%  '
$ws.Range("E13").Value = '  -0.50%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.46'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.31%  '
$ws.Range("D15").Value = '67.559.95'
$ws.Range("E15").Value = '  -0.87%  '
$ws.Range("E16").Value = '  -1.03%  '
$ws.Range("D17").Value = '3.312.14'
$ws.Range("E17").Value = '  +1.23%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '444.15'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +7.80%  '
$ws.Range("E19").Value = '  -0.05%  '
$ws.Range("E20").Value = '  +0.55%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.73'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.24%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.82'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.49%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.12%  '
$ws.Range("D24").Value = '3.455.86'
$ws.Range("E24").Value = '  +1.23%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.513'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.56%  '
$ws.Range("E26").Value = '  +0.67%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.189'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.58%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.05'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.79%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.998'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.28%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.98'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.14%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '22.87'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.54%  '
$ws.Range("E32").Value = '  -2.60%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.998'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.08%  '
$ws.Range("B34").Value = 'Aptos'
$ws.Range("C34").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.81'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.66%  '
$ws.Range("B35").Value = 'Fetch.AI'
$ws.Range("C35").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.24'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.03%  '
$ws.Range("E36").Value = '  +4.44%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '162.86'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.80%  '
$ws.Range("B38").Value = 'Stacks'
$ws.Range("C38").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.85'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.50%  '
$ws.Range("B39").Value = 'EnergySwap'
$ws.Range("C39").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '27.29'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.34%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.793'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.08%  '
$ws.Range("D41").Value = '2.778.26'
$ws.Range("E41").Value = '  +4.36%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.46'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.19%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.30'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.83%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '24.91'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.69%  '
$ws.Range("B45").Value = 'dogwifhat'
$ws.Range("C45").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.41'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.47%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0672'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.11%  '
$ws.Range("B47").Value = 'OKB'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '40.11'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.79%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '328.31'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.50%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0275'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.06%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.993'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.91%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.24'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.32%  '
